$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.317.79"
$ws.Range("E2").Value = "  +2.11%  "
$ws.Range("D3").Value = "2.536.46"
$ws.Range("E3").Value = "  +2.33%  "
$ws.Range("E4").Value = "  +0.26%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "541.10"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.52%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.24"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.11%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.995"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("E8").Value = "  +0.47%  "
$ws.Range("D9").Value = "2.567.09"
$ws.Range("E9").Value = "  +2.69%  "
$ws.Range("E10").Value = "  +1.09%  "
$ws.Range("E11").Value = "  +1.60%  "
$ws.Range("E12").Value = "  -1.72%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.363"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.13%  "
$ws.Range("D14").Value = "2.986.01"
$ws.Range("E14").Value = "  +2.67%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "24.24"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.90%  "
$ws.Range("D16").Value = "60.272.61"
$ws.Range("E16").Value = "  +2.43%  "
$ws.Range("E17").Value = "  +3.12%  "
$ws.Range("D18").Value = "2.538.92"
$ws.Range("E18").Value = "  +2.11%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.36"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.69%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.36"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.33%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "328.68"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.90%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.36%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.93"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.53%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "62.90"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.49%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.440"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.75%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.168"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.64%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.992"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.33%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.04"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.76%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.21"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.50%  "
$ws.Range("D30").Value = "0.0₃0799"
$ws.Range("E30").Value = "  +2.28%  "
$ws.Range("E31").Value = "  +0.70%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.22"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.08%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "162.88"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.80%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.49"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.43%  "
$ws.Range("E35").Value = "  +0.11%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.83"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.26%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.52"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.24%  "
$ws.Range("E38").Value = "  -0.54%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.71"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.06%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "37.19"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.30%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "305.20"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.10%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.843"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.30%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.75"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.48%  "
$ws.Range("E44").Value = "  +2.49%  "
$ws.Range("E45").Value = "  -0.23%  "
$ws.Range("E46").Value = "  +1.20%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "19.08"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.31%  "
$ws.Range("E48").Value = "  +0.94%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "125.10"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.09%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0525"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.11%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0230"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.13%  "
